# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Update DAMSLTag (column I) and DialogAct (column J) for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 5;   Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 8;   Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 9;   Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 27;  Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 29;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 58;  Tag = "aa"; Label = "Agree/Accept" },
    @{ Row = 66;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 70;  Tag = "b";  Label = "Acknowledge (Backchannel)" },
    @{ Row = 78;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 84;  Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 87;  Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 95;  Tag = "b";  Label = "Acknowledge (Backchannel)" },
    @{ Row = 98;  Tag = "%";  Label = "Uninterpretable" },
    @{ Row = 101; Tag = "sd"; Label = "Statement-non-opinion" },
    @{ Row = 104; Tag = "sv"; Label = "Statement-opinion" },
    @{ Row = 116; Tag = "b";  Label = "Acknowledge (Backchannel)" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Label
}

$wb.Save()
